# The sheet gets one new row inserted right before the old row 34, pushing
# the existing "搜狐" .. "大夸科技" block (old rows 34-39) down to 35-40, and
# that new row 34 holds a brand-new shared string "腾讯".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 34 (existing rows 34-39 shift to 35-40).
$ws.Rows.Item(34).Insert() | Out-Null

# Excel's row insert copies the formatting of the row above (which is
# shaded/yellow, style index 5) onto the new row; the target row has no
# explicit style, so strip the inherited formatting back to the default.
$ws.Range("A34").ClearFormats() | Out-Null

# Populate the newly inserted cell.
$ws.Range("A34").Value = "腾讯"

# Match the final selection left in the sheet view.
$ws.Range("C35").Select() | Out-Null
